$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.019246339797974
$ws.Range("B1").Value = 2.401487588882446
$ws.Range("C1").Value = 2.489301204681396
$ws.Range("D1").Value = 3.192335367202759
$ws.Range("E1").Value = 1.321567058563232
